$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 3: new data (G3:K3, M3) mirroring the layout used in row 2
# ---------------------------------------------------------------------------
$ws.Range("G3").Value = 2912000
$ws.Range("H3").Value = 14
$ws.Range("I3").Value = 10
$ws.Range("J3").Formula = "=G3/H3/I3"
$ws.Range("J3").Style = "Berechnung"
$ws.Range("K3").Value = "c8"
$ws.Range("M3").Formula = "=364/H3"
$ws.Range("M3").Style = "Berechnung"

# Apply the "Berechnung" style to M1 / M2 (previously unstyled)
$ws.Range("M1").Style = "Berechnung"
$ws.Range("M2").Style = "Berechnung"

# ---------------------------------------------------------------------------
# New row 11 - continues the A*B*C = D pattern
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = 20000
$ws.Range("B11").Value = 14
$ws.Range("C11").Value = 1
$ws.Range("D11").Formula = "=A11*B11*C11"

# ---------------------------------------------------------------------------
# New row 12 - continues the A*B*C = D pattern
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = 100000
$ws.Range("B12").Value = 14
$ws.Range("C12").Value = 1
$ws.Range("D12").Formula = "=A12*B12*C12"

# ---------------------------------------------------------------------------
# New cell style ("Berechnung" but only left/right borders) for D11 & D12
# ---------------------------------------------------------------------------
$d11 = $ws.Range("D11")
$d11.Style = "Berechnung"
$d11.Borders.Item(8).LineStyle = -4142
$d11.Borders.Item(9).LineStyle = -4142

$d11.Copy()
$ws.Range("D12").PasteSpecial(-4122)

# Restore D12's formula (PasteSpecial only copied formats, not the formula)
$ws.Range("D12").Formula = "=A12*B12*C12"

# ---------------------------------------------------------------------------
# Final selection / active cell
# ---------------------------------------------------------------------------
$ws.Range("E12").Select()
